# Auto-generated edit script applying cryptos.xlsx diff (83 cell updates)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.385.07"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "2.590.40"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'571.72"
$ws.Range("E5").Value = "  +3.06%  "
$ws.Range("D6").Value = "'143.95"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "2.600.90"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("E11").Value = "  +3.54%  "
$ws.Range("E12").Value = "  +11.24%  "
$ws.Range("E13").Value = "  +2.65%  "
$ws.Range("D14").Value = "3.048.66"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").Value = "59.362.97"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").Value = "'22.55"
$ws.Range("E16").Value = "  +7.72%  "
$ws.Range("E17").Value = "  +3.87%  "
$ws.Range("D18").Value = "2.592.79"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("D20").Value = "'335.80"
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "'64.39"
$ws.Range("E24").Value = "  -3.50%  "
$ws.Range("E25").Value = "  +6.87%  "
$ws.Range("D26").Value = "'0.993"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("D29").Value = ("0.0{0}0784" -f [string][char]8323)
$ws.Range("E29").Value = "  +3.39%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "'1.69"
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("D32").Value = "'6.10"
$ws.Range("E32").Value = "  +1.96%  "
$ws.Range("D33").Value = "'158.61"
$ws.Range("E33").Value = "  +2.91%  "
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("E35").Value = "  +3.09%  "
$ws.Range("D36").Value = "'1.16"
$ws.Range("E36").Value = "  +1.61%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "'0.883"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").Value = "'0.876"
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("D39").Value = "'37.19"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("E40").Value = "  +2.32%  "
$ws.Range("D41").Value = "'295.66"
$ws.Range("E41").Value = "  +4.46%  "
$ws.Range("E42").Value = "  +2.16%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").Value = "'0.0981"
$ws.Range("E44").Value = "  +2.79%  "
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("D46").Value = "'0.0539"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("D47").Value = "'19.35"
$ws.Range("E47").Value = "  +2.58%  "
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("E49").Value = "  +6.70%  "
$ws.Range("E50").Value = "  +2.15%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.956.01"
$ws.Range("E51").Value = "  +0.22%  "
